$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6 hold the "Added XML properties" token lists in column G.
# Remove the stray "property" token (a non-tag entry) that was
# mistakenly included right after "lpwstr" in each list.
foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 7)  # column G
    $value = $cell.Value2
    if ($value -ne $null -and $value -like "*,property,*") {
        $cell.Value = $value.Replace(",property,", ",")
    }
}
